# Apply the cryptos.xlsx data refresh: update Price (D) and Volume(1h) (E)
# columns for all data rows, and update Coin (B) / Link (C) for rows whose
# coin ordering changed, per the scraped GitHub Actions update.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{Row=2; B="Bitcoin"; C="https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"; D="29.113.33"; E="  -2.88%  "},
    @{Row=3; B="Ethereum"; C="https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"; D="1.844.39"; E="  -2.14%  "},
    @{Row=4; B="TetherUSD"; C="https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"; D="0.9997"; E="  +0.01%  "},
    @{Row=5; B="XRP"; C="https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"; D="0.6989"; E="  -6.13%  "},
    @{Row=6; B="BNB"; C="https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"; D="237.57"; E="  -2.30%  "},
    @{Row=7; B="USDC"; C="https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"; D="1.000"; E="  +0.00%  "},
    @{Row=8; B="Cardano"; C="https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"; D="0.3030"; E="  -4.24%  "},
    @{Row=9; B="Dogecoin"; C="https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"; D="0.07433"; E="  +2.74%  "},
    @{Row=10; B="Solana"; C="https://coinranking.com/coin/zNZHO_Sjf+solana-sol"; D="23.27"; E="  -6.79%  "},
    @{Row=11; B="TRON"; C="https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"; D="0.08102"; E="  -2.93%  "},
    @{Row=12; B="WrappedEther"; C="https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"; D="1.841.21"; E="  -5.75%  "},
    @{Row=13; B="Polygon"; C="https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"; D="0.7229"; E="  -4.51%  "},
    @{Row=14; B="Polkadot"; C="https://coinranking.com/coin/25W7FG7om+polkadot-dot"; D="5.223"; E="  -3.65%  "},
    @{Row=15; B="Litecoin"; C="https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"; D="88.96"; E="  -4.07%  "},
    @{Row=16; B="WrappedBTC"; C="https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"; D="29.093.57"; E="  -3.05%  "},
    @{Row=17; B="Uniswap"; C="https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"; D="5.792"; E="  -6.41%  "},
    @{Row=18; B="BitcoinCash"; C="https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"; D="240.75"; E="  -3.97%  "},
    @{Row=19; B="ShibaInu"; C="https://coinranking.com/coin/xz24e0BjL+shibainu-shib"; D="0.000007647"; E="  -2.80%  "},
    @{Row=20; B="Avalanche"; C="https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"; D="13.01"; E="  -4.51%  "},
    @{Row=21; B="Dai"; C="https://coinranking.com/coin/MoTuySvg7+dai-dai"; D="0.9995"; E="  -0.07%  "},
    @{Row=22; B="WrappedliquidstakedEther2.0"; C="https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"; D="2.085.00"; E="  -4.40%  "},
    @{Row=23; B="BinanceUSD"; C="https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"; D="1.001"; E="  +0.16%  "},
    @{Row=24; B="Chainlink"; C="https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"; D="7.552"; E="  -6.30%  "},
    @{Row=25; B="Monero"; C="https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"; D="161.98"; E="  -2.25%  "},
    @{Row=26; B="Stellar"; C="https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"; D="0.1463"; E="  -6.58%  "},
    @{Row=27; B="Cosmos"; C="https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D="8.947"; E="  -4.01%  "},
    @{Row=28; B="EthereumClassic"; C="https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"; D="17.98"; E="  -4.09%  "},
    @{Row=29; B="LidoDAOToken"; C="https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"; D="1.934"; E="  -5.39%  "},
    @{Row=30; B="Toncoin"; C="https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D="1.369"; E="  -8.10%  "},
    @{Row=31; B="Filecoin"; C="https://coinranking.com/coin/ymQub4fuB+filecoin-fil"; D="4.474"; E="  -3.28%  "},
    @{Row=32; B="PancakeSwap"; C="https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"; D="1.491"; E="  -3.03%  "},
    @{Row=33; B="InternetComputer(DFINITY)"; C="https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"; D="4.018"; E="  -5.21%  "},
    @{Row=34; B="Hedera"; C="https://coinranking.com/coin/jad286TjB+hedera-hbar"; D="0.05175"; E="  -3.75%  "},
    @{Row=35; B="ARBITRUM"; C="https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"; D="1.183"; E="  -5.71%  "},
    @{Row=36; B="ImmutableX"; C="https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"; D="0.7092"; E="  -6.59%  "},
    @{Row=37; B="Frax"; C="https://coinranking.com/coin/KfWtaeV1W+frax-frax"; D="1.000"; E="  +0.10%  "},
    @{Row=38; B="HuobiToken"; C="https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"; D="2.648"; E="  -2.24%  "},
    @{Row=39; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.01869"; E="  -5.23%  "},
    @{Row=40; B="MXToken"; C="https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"; D="2.675"; E="  -3.18%  "},
    @{Row=41; B="TrustWalletToken"; C="https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"; D="0.9012"; E="  +4.11%  "},
    @{Row=42; B="TheSandbox"; C="https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"; D="0.4281"; E="  -6.42%  "},
    @{Row=43; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="5.921"; E="  -2.68%  "},
    @{Row=44; B="Maker"; C="https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"; D="1.054.57"; E="  -4.68%  "},
    @{Row=45; B="Aave"; C="https://coinranking.com/coin/ixgUfzmLR+aave-aave"; D="69.80"; E="  -4.36%  "},
    @{Row=46; B="PaxDollar"; C="https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"; D="0.9998"; E="  -0.07%  "},
    @{Row=47; B="Quant"; C="https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"; D="101.51"; E="  -2.99%  "},
    @{Row=48; B="RenderToken"; C="https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"; D="1.751"; E="  -6.42%  "},
    @{Row=49; B="EnergySwap"; C="https://coinranking.com/coin/SbWqqTui-+energyswap-ens"; D="9.240"; E="  -3.17%  "},
    @{Row=50; B="Aptos"; C="https://coinranking.com/coin/HGYj5JCv5+aptos-apt"; D="7.056"; E="  -7.47%  "},
    @{Row=51; B="RocketPoolETH"; C="https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"; D="1.975.79"; E="  -4.85%  "}
)

foreach ($item in $rows) {
    $r = $item.Row

    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C

    # Column D ("Price") often looks numeric (e.g. "1.000", "0.9997",
    # "29.113.33"); force text format first so Excel keeps the exact
    # literal string instead of re-interpreting it as a number/date.
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $item.D

    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $item.E
}
